$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 73077
$ws.Range("E2").Value = 8747
$ws.Range("F2").Value = 12

# Row 3
$ws.Range("C3").Value = 58397
$ws.Range("E3").Value = 3408
$ws.Range("F3").Value = 5.8

# Row 4
$ws.Range("C4").Value = 451779
$ws.Range("E4").Value = 243979
$ws.Range("F4").Value = 54

# Row 5
$ws.Range("C5").Value = 78642
$ws.Range("E5").Value = 25949
$ws.Range("F5").Value = 33

# Row 6
$ws.Range("C6").Value = 50473
$ws.Range("E6").Value = 24561
$ws.Range("F6").Value = 48.7

# Row 7
$ws.Range("C7").Value = 282888
$ws.Range("E7").Value = 10614
$ws.Range("F7").Value = 3.8

# Row 8
$ws.Range("C8").Value = 70064
$ws.Range("E8").Value = 3068
$ws.Range("F8").Value = 4.4

# Row 9
$ws.Range("C9").Value = 52280
$ws.Range("E9").Value = -1214
$ws.Range("F9").Value = -2.3

# Row 10
$ws.Range("C10").Value = 425484
$ws.Range("D10").Value = 494604
$ws.Range("E10").Value = 69120
$ws.Range("F10").Value = 16.2

# Row 11
$ws.Range("C11").Value = 88975
$ws.Range("D11").Value = 95411
$ws.Range("E11").Value = 6436
$ws.Range("F11").Value = 7.2

# Row 12
$ws.Range("C12").Value = 67893
$ws.Range("D12").Value = 75971
$ws.Range("E12").Value = 8078
$ws.Range("F12").Value = 11.9
